# Add non-resident ("2020_non_res") abortion-by-state data as a new worksheet
# at the end of the workbook, matching the reports update that accompanies
# the 2019-2021 non-resident data added to the script/csv/index.html.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the current last tab so it lands at the end of
# the workbook (a plain Add() would insert before the active sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2020_non_res"

# Header row
$ws.Range("A1").Value = "state"
$ws.Range("B1").Value = "abortions"

# Non-resident abortion counts by state of residence
$ws.Range("A2").Value = "Kentucky"
$ws.Range("B2").Value = 200

$ws.Range("A3").Value = "Illinois"
$ws.Range("B3").Value = 75

$ws.Range("A4").Value = "Michigan"
$ws.Range("B4").Value = 55

$ws.Range("A5").Value = "Ohio"
$ws.Range("B5").Value = 27

$ws.Range("A6").Value = "Tennessee"
$ws.Range("B6").Value = 16

$ws.Range("A7").Value = "Other"
$ws.Range("B7").Value = 11

# Leave the selection the way Excel would after typing the last value and
# hitting enter once more, with the new sheet active/selected.
$ws.Range("B8").Select()
